$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "struct_list": selection moves from F16 to the whole row 5
# (do this before activating "ai" -- selecting on a sheet makes it
# the active sheet, and "ai" must end up as the final active tab).
# ------------------------------------------------------------------
$wsStruct = $wb.Worksheets.Item("struct_list")
$wsStruct.Rows.Item(5).Select() | Out-Null

# ------------------------------------------------------------------
# Sheet "ai": pseudo buffering algorithm rows (8-11)
# Cell-by-cell order matters: it reproduces the exact shared-string
# table append order seen in the target workbook.
# ------------------------------------------------------------------
$wsAi = $wb.Worksheets.Item("ai")

function Set-TextValue {
    param($range, [string]$text)
    # Route through a text formula + paste-as-values so Excel stores a
    # literal text shared string (t="s") instead of auto-coercing
    # look-alike words ("true"/"false") into a Boolean cell, while
    # keeping the cell's existing style untouched.
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $wb.Application.CutCopyMode = $false
}

# Row 11
$wsAi.Range("A11").Value = "meas_calib_b"
# Row 10
$wsAi.Range("A10").Value = "meas_calib_a"
# Row 9
$wsAi.Range("A9").Value = "meas_value"
# Descriptions (column B), rows 9, 10, 11
$wsAi.Range("B9").Value = "Engineer units value"
$wsAi.Range("B10").Value = "Engineer units add"
$wsAi.Range("B11").Value = "Engineer units mul"
# Default value for meas_calib_b
$wsAi.Range("G11").Value = "{1.0f,1.0f}"
# Row 8 (avg_num), added last
$wsAi.Range("A8").Value = "avg_num"

# Remaining cells, reusing already-existing shared strings / numbers.
$wsAi.Range("C8").Value = "u16"
$wsAi.Range("D8").Value = 2
Set-TextValue $wsAi.Range("E8") "false"
$wsAi.Range("F8").Value = 1018
$wsAi.Range("G8").Value = 10
$wsAi.Range("H8").Value = 1
$wsAi.Range("I8").Value = 100
$wsAi.Range("J8").Value = "auto"

$wsAi.Range("C9").Value = "float"
$wsAi.Range("D9").Value = 2
Set-TextValue $wsAi.Range("E9") "true"
$wsAi.Range("F9").Value = 1020

$wsAi.Range("C10").Value = "float"
$wsAi.Range("D10").Value = 2
Set-TextValue $wsAi.Range("E10") "false"
$wsAi.Range("F10").Value = 1024
$wsAi.Range("G10").Value = "{0.0f,0.0f}"
$wsAi.Range("J10").Value = "auto"

$wsAi.Range("C11").Value = "float"
$wsAi.Range("D11").Value = 2
Set-TextValue $wsAi.Range("E11") "false"
$wsAi.Range("F11").Value = 1028
$wsAi.Range("J11").Value = "auto"

# ------------------------------------------------------------------
# Make "ai" the active sheet/tab and set its new selection (F12)
# ------------------------------------------------------------------
$wsAi.Activate()
$wsAi.Range("F12").Select() | Out-Null
